$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values per row (rows 2-25, columns B-N), applied for the 380 kV case
$row = @(5.689121947306627,0,4.668847530398286,11.7275782455013,25.86160637963263,31.99044886572986,14.67768628747004,19.85579666574063,0,10.7472741761131,0,14.36069988942034,18.35393530346338)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B2:N2").Value2 = $arr

$row = @(5.593870226478985,0,4.662562700883706,11.58956654610724,25.73412087334159,31.70567781638372,14.69607122223506,19.94310023892019,0,10.35765294126347,0,14.11901570511218,18.41906948097751)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B3:N3").Value2 = $arr

$row = @(5.535534133082446,0,4.658672058454461,11.50803411085259,25.66391029850729,31.54193380342438,14.7106387114682,19.99991304736397,0,10.11316290774348,0,13.9724590500224,18.46085650972473)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B4:N4").Value2 = $arr

$row = @(5.511830591205262,0,4.657078994076895,11.4756536512131,25.63734914224783,31.47807028646255,14.71739780196716,20.02387198727148,0,10.01237354392683,0,13.9132798991922,18.47833773062245)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B5:N5").Value2 = $arr

$row = @(5.507899722628264,0,4.656814019429432,11.47032899047472,25.63306308450017,31.46764067603923,14.71856977530675,20.02789912471066,0,9.995572755016369,0,13.90348847067883,18.4812678568047)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B6:N6").Value2 = $arr

$row = @(5.535214138075757,0,4.658650604095692,11.50759394774716,25.66354375783671,31.54106083535282,14.71072653842593,20.00023289669209,0,10.111808084967,0,13.97165862969677,18.461090432843)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B7:N7").Value2 = $arr

$row = @(5.656267017338265,0,4.666687334747543,11.67934965468036,25.81598948319034,31.88999910985631,14.68334382249142,19.88523344231325,0,10.61411908187651,0,14.27703464819792,18.37602202537952)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B8:N8").Value2 = $arr

$row = @(5.893473519498703,0,4.682179388556872,12.03983988054727,26.17780928583074,32.65896336470323,14.65573563846343,19.68514835786852,0,11.55092383506484,0,14.88687962474501,18.22337298136812)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B9:N9").Value2 = $arr

$row = @(6.065869820460742,0,4.693377298103464,12.31665214477325,26.48023800982464,33.27055834705266,14.65143865244714,19.55361112977743,0,12.20217302347636,0,15.33689966377698,18.11976330687041)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B10:N10").Value2 = $arr

$row = @(6.1435329468595,0,4.698427224129411,12.44463334450542,26.62533699991516,33.55775613476016,14.65296657295528,19.49712251815085,0,12.48907936879153,0,15.54108128152442,18.0744625833973)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B11:N11").Value2 = $arr

$row = @(6.172804422225886,0,4.7003327872747,12.49334635304382,26.68132371518984,33.66769986173628,14.65404629689216,19.47621287462731,0,12.59628593794214,0,15.61824507872954,18.05757025728943)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B12:N12").Value2 = $arr

$row = @(6.166506856578892,0,4.69992269597304,12.48284478713084,26.66922035613067,33.64397033442334,14.65379146923776,19.48069473292438,0,12.57326244133059,0,15.60163459861576,18.06119668736984)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B13:N13").Value2 = $arr

$row = @(6.145944049976324,0,4.698584133193432,12.44863626288099,26.62992244222727,33.56677805743289,14.65304535831954,19.49539262238724,0,12.49792864698054,0,15.54743314500683,18.07306759684598)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B14:N14").Value2 = $arr

$row = @(6.133329961491093,0,4.697763336065113,12.42771364062008,26.60598563571313,33.51964719305111,14.65265361037246,19.50445818961295,0,12.45159458029741,0,15.5142106844249,18.08037296187536)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B15:N15").Value2 = $arr

$row = @(6.060776579474256,0,4.693046352834157,12.30832569449133,26.47090338441557,33.25196123212417,14.65140891513617,19.55737015926431,0,12.18322655162271,0,15.32353837635912,18.12276056207021)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B16:N16").Value2 = $arr

$row = @(6.016051895101667,0,4.690141067245117,12.23557893628589,26.38993320661647,33.08997194590903,14.65153766735737,19.59068742879722,0,12.01612585361525,0,15.20637332863823,18.14923220211351)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B17:N17").Value2 = $arr

$row = @(5.990257448752109,0,4.688465900341565,12.19393341905016,26.34407197403657,32.9976515989826,14.65193952598196,19.61016580618295,0,11.91913869255951,0,15.13893726661283,18.16463049000443)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B18:N18").Value2 = $arr

$row = @(5.981512760470913,0,4.687898021613441,12.17986813292297,26.32866741387909,32.96654288559467,14.65213187337684,19.61681498437709,0,11.88615324468033,0,15.11609914554361,18.16987375391352)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B19:N19").Value2 = $arr

$row = @(6.020820368324869,0,4.690450769501004,12.24330295470014,26.39847935693401,33.10712858413204,14.65149003147336,19.58710813114278,0,12.03400531380237,0,15.21885106903227,18.14639640738332)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B20:N20").Value2 = $arr

$row = @(6.15198780644067,0,4.698977487454877,12.45867774011675,26.64143728014407,33.58941988688441,14.65325090757015,19.49106243250752,0,12.52009574376236,0,15.56335825132186,18.06957372098756)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B21:N21").Value2 = $arr

$row = @(6.236897498319681,0,4.704510755520104,12.60086934951589,26.80627083633627,33.91150145383693,14.65732264757355,19.43109677292906,0,12.88588548093828,0,15.78757404981585,18.02089261053564)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B22:N22").Value2 = $arr

$row = @(6.191663293049769,0,4.701561285428933,12.5248633402296,26.71775689148409,33.73900587015875,14.65488219744074,19.46284487454895,0,12.68003620617124,0,15.66801690412673,18.04673535961127)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B23:N23").Value2 = $arr

$row = @(6.018664793674419,0,4.690310768249845,12.2398103693585,26.39461349113525,33.09936953866385,14.65151054642266,19.5887253228345,0,12.02592486424899,0,15.21321011427832,18.14767791124329)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B24:N24").Value2 = $arr

$row = @(5.829508202946831,0,4.678019217011641,11.94003405645526,26.07336272300815,32.44236694036397,14.66040166636359,19.73655832262305,0,11.30349980859693,0,14.72121757378609,18.26316164651983)
$arr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $arr[0,$c] = $row[$c] }
$ws.Range("B25:N25").Value2 = $arr

